# Insert a new data row for Zanahoria / Terminal La Palmera de La Serena
# at sheet row 224, shifting all subsequent rows (224-355) down by one
# (they become rows 225-356). The new row carries a new weekly reading;
# the previously-last row (355) is duplicated down into the new row 356
# as a consequence of the shift, exactly matching the existing data below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 224 (and everything below it) down by one row.
$ws.Rows.Item(224).Insert()

# Populate the newly-blank row 224 with the new record.
$ws.Range("A224").Value = 8
$ws.Range("B224").Value = "Terminal La Palmera de La Serena"
$ws.Range("C224").Value = "Coquimbo"
$ws.Range("D224").Value = 44719
$ws.Range("E224").Value = 4
$ws.Range("F224").Value = 100114013
$ws.Range("G224").Value = "Zanahoria"
$ws.Range("H224").Value = "Sin especificar"
$ws.Range("I224").Value = "Primera"
$ws.Range("J224").Value = 500
$ws.Range("K224").Value = 6000
$ws.Range("L224").Value = 7000
$ws.Range("M224").Value = 6500
$ws.Range("N224").Value = "`$/saco 20 kilos"
$ws.Range("O224").Value = "Provincia del Elquí"
$ws.Range("P224").Value = 325
$ws.Range("Q224").Value = 20
$ws.Range("R224").Value = "Hortaliza"
